$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 39892.5
$ws.Range("J3").Value = 39892.5
$ws.Range("L3").Value = 39892.5
$ws.Range("N3").Value = -40120.5
# Row 6
$ws.Range("H6").Value = 95.29412000000001
$ws.Range("I6").Value = 108.545456
$ws.Range("J6").Value = 71
$ws.Range("K6").Value = 325.636368
$ws.Range("L6").Value = 213
$ws.Range("M6").Value = -213.636368
$ws.Range("N6").Value = -437
# Row 28
$ws.Range("H28").Value = 10150.3
$ws.Range("I28").Value = 1751.25
$ws.Range("J28").Value = 15749.667
$ws.Range("K28").Value = 1751.25
$ws.Range("L28").Value = 15749.667
$ws.Range("M28").Value = -1266.25
$ws.Range("N28").Value = -16719.667
# Row 74
$ws.Range("H74").Value = 7560.2
$ws.Range("I74").Value = 7560.2
$ws.Range("K74").Value = 7560.2
$ws.Range("M74").Value = -6624.2
# Row 77
$ws.Range("H77").Value = 7560.2
$ws.Range("I77").Value = 7560.2
$ws.Range("K77").Value = 37801
$ws.Range("M77").Value = -33121
# Row 97
$ws.Range("H97").Value = 749.8
$ws.Range("J97").Value = 749.8
$ws.Range("L97").Value = 2249.4
$ws.Range("N97").Value = -3241.4
# Row 102
$ws.Range("H102").Value = 39892.5
$ws.Range("J102").Value = 39892.5
$ws.Range("L102").Value = 39892.5
$ws.Range("N102").Value = -46382.5
# Row 107
$ws.Range("H107").Value = 508.25
$ws.Range("I107").Value = 477.7143
$ws.Range("K107").Value = 477.7143
$ws.Range("M107").Value = 1442.2857
# Row 141
$ws.Range("H141").Value = 2841.3333
$ws.Range("I141").Value = 1794.25
$ws.Range("K141").Value = 5382.75
$ws.Range("M141").Value = -202.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1132
$ws.Range("I2").Value = 1132
$ws.Range("K2").Value = 1132
$ws.Range("M2").Value = -1019
# Row 45
$ws.Range("H45").Value = 2694.2856
$ws.Range("I45").Value = 2376
$ws.Range("J45").Value = 3490
$ws.Range("K45").Value = 2376
$ws.Range("L45").Value = 3490
$ws.Range("M45").Value = -1999
$ws.Range("N45").Value = -4244
# Row 61
$ws.Range("H61").Value = 2579.8333
$ws.Range("I61").Value = 2642.2727
$ws.Range("K61").Value = 2642.2727
$ws.Range("M61").Value = -2430.2727
# Row 97
$ws.Range("H97").Value = 980
$ws.Range("I97").Value = 857.61536
$ws.Range("K97").Value = 857.61536
$ws.Range("M97").Value = -361.61536
# Row 102
$ws.Range("H102").Value = 9002
$ws.Range("I102").Value = 5010
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 5010
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -3388
$ws.Range("N102").Value = -13244
# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = ""
$ws.Range("N105").Value = 0
# Row 116
$ws.Range("H116").Value = 1132
$ws.Range("I116").Value = 1132
$ws.Range("K116").Value = 1132
$ws.Range("M116").Value = 1162
# Row 122
$ws.Range("H122").Value = 1907.3572
$ws.Range("I122").Value = 1580.8182
$ws.Range("K122").Value = 4742.4546
$ws.Range("M122").Value = -2292.4546
# Row 132
$ws.Range("H132").Value = 2592
$ws.Range("I132").Value = 2592
$ws.Range("K132").Value = 7776
$ws.Range("M132").Value = -5246
# Row 136
$ws.Range("H136").Value = 2579.8333
$ws.Range("I136").Value = 2642.2727
$ws.Range("K136").Value = 7926.8181
$ws.Range("M136").Value = -5376.8181

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1132
$ws.Range("I3").Value = 1132
$ws.Range("K3").Value = 1132
$ws.Range("M3").Value = -1018
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").Value = ""
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").Value = ""
# Row 86
$ws.Range("H86").Value = 8333.5
$ws.Range("I86").Value = 5000.5
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 5000.5
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -3877.5
$ws.Range("N86").Value = -12246
# Row 89
$ws.Range("H89").Value = 8333.5
$ws.Range("I89").Value = 5000.5
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 25002.5
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -19386.5
$ws.Range("N89").Value = -61232
# Row 94
$ws.Range("H94").Value = 380.85715
$ws.Range("I94").Value = 371.69232
$ws.Range("J94").Value = 500
$ws.Range("K94").Value = 371.69232
$ws.Range("L94").Value = 500
$ws.Range("M94").Value = 79.30768
$ws.Range("N94").Value = -1402

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6423.9414
$ws.Range("I31").Value = 3826.75
$ws.Range("K31").Value = 3826.75
$ws.Range("M31").Value = -3531.75
# Row 34
$ws.Range("H34").Value = 6423.9414
$ws.Range("I34").Value = 3826.75
$ws.Range("K34").Value = 3826.75
$ws.Range("M34").Value = -3624.75
# Row 43
$ws.Range("H43").Value = 17597.834
$ws.Range("J43").Value = 17597.834
$ws.Range("L43").Value = 17597.834
$ws.Range("N43").Value = -17965.834
# Row 51
$ws.Range("H51").Value = 31171
$ws.Range("I51").Value = 1090
$ws.Range("J51").Value = 38691.25
$ws.Range("K51").Value = 1090
$ws.Range("L51").Value = 38691.25
$ws.Range("M51").Value = -354
$ws.Range("N51").Value = -40163.25
# Row 61
$ws.Range("H61").Value = 31171
$ws.Range("I61").Value = 1090
$ws.Range("J61").Value = 38691.25
$ws.Range("K61").Value = 1090
$ws.Range("L61").Value = 38691.25
$ws.Range("M61").Value = -742
$ws.Range("N61").Value = -39387.25
# Row 101
$ws.Range("H101").Value = 17597.834
$ws.Range("J101").Value = 17597.834
$ws.Range("L101").Value = 17597.834
$ws.Range("N101").Value = -24087.834
# Row 107
$ws.Range("H107").Value = 532.8570999999999
$ws.Range("I107").Value = 264
$ws.Range("J107").Value = 734.5
$ws.Range("K107").Value = 264
$ws.Range("L107").Value = 734.5
$ws.Range("M107").Value = 1656
$ws.Range("N107").Value = -4574.5
# Row 132
$ws.Range("H132").Value = 1447.3334
$ws.Range("I132").Value = 1447.3334
$ws.Range("K132").Value = 4342.0002
$ws.Range("M132").Value = -1812.0002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 54
$ws.Range("H54").Value = 2887.5
# Row 60
$ws.Range("H60").Value = 985.6667
$ws.Range("I60").Value = 237.2
$ws.Range("J60").Value = 2233.111
$ws.Range("K60").Value = 711.5999999999999
$ws.Range("L60").Value = 6699.333
$ws.Range("M60").Value = -460.5999999999999
$ws.Range("N60").Value = -7201.333
# Row 107
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = ""
$ws.Range("N107").Value = 0

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 9834.666999999999
$ws.Range("I113").Value = 9669.333000000001
$ws.Range("K113").Value = 9669.333000000001
$ws.Range("M113").Value = -7499.333000000001
# Row 132
$ws.Range("H132").Value = 59369.832
$ws.Range("I132").Value = 69977.13
$ws.Range("K132").Value = 209931.39
$ws.Range("M132").Value = -207401.39

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 4363.6665
$ws.Range("I16").Value = 1545
$ws.Range("J16").Value = 10001
$ws.Range("K16").Value = 1545
$ws.Range("L16").Value = 10001
$ws.Range("M16").Value = -1375
$ws.Range("N16").Value = -10341
# Row 44
$ws.Range("H44").Value = 7000
$ws.Range("J44").Value = 7000
$ws.Range("L44").Value = 7000
$ws.Range("N44").Value = -7912
# Row 55
$ws.Range("H55").Value = 1278.0952
$ws.Range("I55").Value = 1113.5454
$ws.Range("J55").Value = 1459.1
$ws.Range("K55").Value = 1113.5454
$ws.Range("L55").Value = 1459.1
$ws.Range("M55").Value = -940.5454
$ws.Range("N55").Value = -1805.1
# Row 139
$ws.Range("H139").Value = 615000
$ws.Range("J139").Value = 615000
$ws.Range("L139").Value = 615000
$ws.Range("N139").Value = -625280

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = ""
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = ""
# Row 81
$ws.Range("H81").Value = 735.4
$ws.Range("I81").Value = 735.4
$ws.Range("K81").Value = 1470.8
$ws.Range("M81").Value = -409.8
# Row 84
$ws.Range("H84").Value = 735.4
$ws.Range("I84").Value = 735.4
$ws.Range("K84").Value = 7354
$ws.Range("M84").Value = -2050
# Row 107
$ws.Range("H107").Value = 661.7273
$ws.Range("I107").Value = 517.5714
$ws.Range("J107").Value = 914
$ws.Range("K107").Value = 1552.7142
$ws.Range("L107").Value = 2742
$ws.Range("M107").Value = 367.2857999999999
$ws.Range("N107").Value = -6582
# Row 119
$ws.Range("H119").Value = 10698
$ws.Range("J119").Value = 10698
$ws.Range("L119").Value = 10698
$ws.Range("N119").Value = -20374
# Row 136
$ws.Range("H136").Value = 2262.4092
$ws.Range("I136").Value = 1658.8438
$ws.Range("K136").Value = 4976.5314
$ws.Range("M136").Value = -2426.5314
$ws.Range("N136").Value = -16715.7501

